# RPA datasets push 2024-05-08
# Insert a new IPO record (SK증권제12호스팩) into the "01_IB전략컨설팅부" sheet
# as row 3, pushing the existing rows 3-13 down to rows 4-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row above the current row 3 (민테크), shifting it and
# everything below it down by one row.
$ws.Rows.Item(3).Insert()

# --- Populate the new row 3 with the SK증권제12호스팩 data ---

# Date-like text columns: Excel's COM layer auto-parses strings such as
# "2024-04-17" into date serials when assigned directly, so we prefix with
# a quote to force literal text entry (matching the source data, which is
# plain text), then copy the style from an existing unstyled date cell so
# the cell keeps the default "General" style (no stray number format).
$ws.Cells.Item(3,1).Value = "'2024-04-17"
$ws.Cells.Item(3,1).Style = $ws.Cells.Item(2,1).Style

$ws.Cells.Item(3,2).Value = "'2024-04-18"
$ws.Cells.Item(3,2).Style = $ws.Cells.Item(2,2).Style

$ws.Cells.Item(3,3).Value = "'2024-05-07"
$ws.Cells.Item(3,3).Style = $ws.Cells.Item(2,3).Style

$ws.Cells.Item(3,4).Value = "SK"
$ws.Cells.Item(3,5).Value = "SK증권제12호스팩"

$ws.Cells.Item(3,6).Value = 3000000
$ws.Cells.Item(3,7).Value = 3000000
$ws.Cells.Item(3,8).Value = 0
$ws.Cells.Item(3,9).Value = 2000
$ws.Cells.Item(3,10).Value = 2000
$ws.Cells.Item(3,11).Value = 3310000
$ws.Cells.Item(3,12).Value = 0
$ws.Cells.Item(3,13).Value = 2000

$ws.Cells.Item(3,14).Value = "1,189.41:1"
$ws.Cells.Item(3,15).Value = "-"

$ws.Cells.Item(3,16).Value = 0
$ws.Cells.Item(3,17).Value = 0
$ws.Cells.Item(3,18).Value = 0
$ws.Cells.Item(3,19).Value = 0
$ws.Cells.Item(3,20).Value = 0
$ws.Cells.Item(3,21).Value = 0
$ws.Cells.Item(3,22).Value = 0
$ws.Cells.Item(3,23).Value = 0
$ws.Cells.Item(3,24).Value = 0

$ws.Cells.Item(3,25).Value = "기업인수목적 주식회사"
